# Increment the "想去人数" (want-to-go count, column F) by 1 for a handful of
# rows on both the "展览" and "全部类型" worksheets, matching the refreshed
# data snapshot published to gh-pages.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> (old value, new value) for column F on each affected sheet.
$updates = @{
    2  = 113
    17 = 715
    19 = 351
    20 = 4361
    26 = 730
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
